# Update the CyberMAGICS participant-introduction template:
#  - the "datetimeFigureOut" date placeholder cached on the slide master
#    and on every slide layout: 6/17/24 -> 6/5/25
#  - the workshop date mentioned in the slide 1 caption textbox:
#    "Workshop, June 25-26, 2024" -> "Workshop, June 5-6, 2025"

$p = $ppt.ActivePresentation

$oldDate = "6/17/24"
$newDate = "6/5/25"

# --- Slide Master: fix the cached date placeholder text ---
foreach ($sh in $p.SlideMaster.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every Slide Layout: fix the cached date placeholder text ---
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    foreach ($sh in $cl.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide 1: update the workshop date in the "CyberMAGICS Workshop, ..." caption ---
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 35")
$tr = $sh.TextFrame.TextRange
$prefix = "CyberMAGICS"
$rest = $tr.Characters($prefix.Length + 1, $tr.Length - $prefix.Length)
$rest.Text = " Workshop, June 5-6, 2025"
